# Generate Report for Handoff
# Unify the "Latest Handoff Date"/"Latest Handoff Datetime" timestamps that
# were previously split into two slightly different values into a single,
# newer timestamp value (simulating a fresh report generation run where all
# of these files now share the same latest-handoff timestamp).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$newDate1 = "2016-29-20 14:29:27"
foreach ($r in @(7, 10, 11, 12, 13, 14, 15, 16)) {
    $ws1.Cells.Item($r, 4).Value = $newDate1
}

# --- zh-cn sheet -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$newDate2 = "2016-03-20 14:29:23"
foreach ($r in @(7, 10, 11, 12, 13, 14, 15, 16)) {
    $ws2.Cells.Item($r, 5).Value = $newDate2
}

# --- de-de sheet -------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$newDate3 = "2016-03-20 14:29:27"
foreach ($r in @(7, 10, 11, 12, 13, 14, 15, 16)) {
    $ws3.Cells.Item($r, 5).Value = $newDate3
}
